# hours update and TAR update
# Adds two new status-report rows (1/26/2010) to Sheet1:
#   Row 18: 1/26/2010 | 2 hours | Group Meeting
#   Row 19: 1/26/2010 | 1 hour  | Weekly Meeting

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The dates in column A on this sheet are stored as literal text (e.g.
# "1/21/2010" on row 17), not as real date serial numbers. Pre-format the
# two new A cells as Text so Excel's automatic date recognition doesn't
# convert the typed string into a date value, then restore the default
# formatting (so the cells keep the workbook's normal/general style, same
# as the existing text-date rows).
$ws.Range("A18:A19").NumberFormat = "@"

$ws.Range("A18").Value = "1/26/2010"
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "Group Meeting"

$ws.Range("A19").Value = "1/26/2010"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = "Weekly Meeting"

$ws.Range("A18:A19").ClearFormats()

# Move the active selection down to A20, matching where entry would
# continue after the two new rows.
$ws.Range("A20").Select()
